$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Column M holds generated-code formulas of the form:
#   ="if ("&K<r>&" == 20) { "&J<r>&".setText("&CHAR(34)&"X"&CHAR(34)&"); } ...
# The strike threshold in the first condition changes from 20 to 10.
# A handful of rows (22,23,45,46,57,58) have no formula in column M and
# must stay empty, so we update every other row individually -- Excel
# then re-forms the shared-formula groups on its own, exactly as before.

$emptyRows = @(22, 23, 45, 46, 57, 58)

for ($row = 1; $row -le 68; $row++) {
    if ($emptyRows -contains $row) { continue }

    $formula = '="if ("&K' + $row + '&" == 10) { "&J' + $row + '&".setText("&CHAR(34)&"X"&CHAR(34)&"); } else if ("&K' + $row + '&" == 30)  { "&J' + $row + '&".setText("&CHAR(34)&"/"&CHAR(34)&"); } else if ("&K' + $row + '&" == 40)  { "&J' + $row + '&".setText("&CHAR(34)&CHAR(34)&"); } else if ("&K' + $row + '&" == 0)  { "&J' + $row + '&".setText("&CHAR(34)&"-"&CHAR(34)&"); }  else {"&J' + $row + '&".setText("&CHAR(34)&CHAR(34)&"+"&K' + $row + '&");}"'

    $ws.Range("M$row").Formula = $formula
}

$excel.Calculate()

# Restore the scroll position of the sheet view (topLeftCell moved up from E40 to E37).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 5
